$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.473.23"
$ws.Range("E2").Value = "  -2.70%  "
$ws.Range("D3").Value = "2.952.65"
$ws.Range("E3").Value = "  -3.69%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'496.07"
$ws.Range("E5").Value = "  -6.10%  "
$ws.Range("D6").Value = "'133.85"
$ws.Range("E6").Value = "  -6.98%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "'0.423"
$ws.Range("E8").Value = "  -5.75%  "
$ws.Range("D9").Value = "'7.12"
$ws.Range("E9").Value = "  -6.87%  "
$ws.Range("D10").Value = "'0.105"
$ws.Range("E10").Value = "  -6.80%  "
$ws.Range("D11").Value = "'0.351"
$ws.Range("E11").Value = "  -5.73%  "
$ws.Range("D12").Value = "3.464.36"
$ws.Range("E12").Value = "  -3.48%  "
$ws.Range("E13").Value = "  -3.46%  "
$ws.Range("D14").Value = "'25.78"
$ws.Range("E14").Value = "  -5.91%  "
$ws.Range("D15").Value = "'0.0000156"
$ws.Range("E15").Value = "  -9.67%  "
$ws.Range("D16").Value = "56.556.04"
$ws.Range("E16").Value = "  -2.47%  "
$ws.Range("D17").Value = "2.958.74"
$ws.Range("E17").Value = "  -3.39%  "
$ws.Range("D18").Value = "'5.95"
$ws.Range("E18").Value = "  -4.24%  "
$ws.Range("D19").Value = "'12.41"
$ws.Range("E19").Value = "  -6.28%  "
$ws.Range("D20").Value = "'7.73"
$ws.Range("E20").Value = "  -5.98%  "
$ws.Range("D21").Value = "'315.42"
$ws.Range("E21").Value = "  -8.05%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "'0.483"
$ws.Range("E24").Value = "  -4.56%  "
$ws.Range("D25").Value = "'62.48"
$ws.Range("E25").Value = "  -4.59%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'0.161"
$ws.Range("E27").Value = "  -5.72%  "
$ws.Range("D28").Value = "0.0₃0856"
$ws.Range("E28").Value = "  -12.88%  "
$ws.Range("D29").Value = "'6.44"
$ws.Range("E29").Value = "  -8.56%  "
$ws.Range("D30").Value = "'6.99"
$ws.Range("E30").Value = "  -6.22%  "
$ws.Range("D31").Value = "'1.74"
$ws.Range("E31").Value = "  -6.38%  "
$ws.Range("D32").Value = "'19.85"
$ws.Range("E32").Value = "  -6.75%  "
$ws.Range("D33").Value = "'1.12"
$ws.Range("E33").Value = "  -9.43%  "
$ws.Range("D34").Value = "'152.64"
$ws.Range("E34").Value = "  -3.31%  "
$ws.Range("D35").Value = "'4.45"
$ws.Range("E35").Value = "  -7.96%  "
$ws.Range("D36").Value = "'5.66"
$ws.Range("E36").Value = "  -5.61%  "
$ws.Range("E37").Value = "  -10.18%  "
$ws.Range("D38").Value = "'23.75"
$ws.Range("E38").Value = "  -9.16%  "
$ws.Range("D39").Value = "'0.0649"
$ws.Range("E39").Value = "  -7.67%  "
$ws.Range("D40").Value = "2.987.81"
$ws.Range("E40").Value = "  -3.61%  "
$ws.Range("D41").Value = "'37.22"
$ws.Range("E41").Value = "  -1.62%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").Value = "'0.638"
$ws.Range("E43").Value = "  -4.45%  "
$ws.Range("D44").Value = "'3.66"
$ws.Range("E44").Value = "  -7.57%  "
$ws.Range("D45").Value = "2.141.70"
$ws.Range("E45").Value = "  -8.82%  "
$ws.Range("D46").Value = "'1.34"
$ws.Range("E46").Value = "  -9.52%  "
$ws.Range("D47").Value = "'5.83"
$ws.Range("E47").Value = "  -4.97%  "
$ws.Range("D48").Value = "'0.915"
$ws.Range("E48").Value = "  -12.79%  "
$ws.Range("D49").Value = "'0.0229"
$ws.Range("E49").Value = "  -6.86%  "
$ws.Range("D50").Value = "'18.87"
$ws.Range("E50").Value = "  -7.05%  "
$ws.Range("D51").Value = "'1.72"
$ws.Range("E51").Value = "  -14.70%  "
